# Commit: "Added getAll to rest"
# Adds a new time-tracking entry (row 44) to the Arbeitszeiten sheet:
#   Datum = 2018-10-18, Name = David, Was = "GetAll hinzugefügt", Dauer = 1.5h

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A44").Value = Get-Date -Year 2018 -Month 10 -Day 18 -Hour 0 -Minute 0 -Second 0
$ws.Range("B44").Value = "David"
$ws.Range("C44").Value = "GetAll hinzugefügt"
$ws.Range("D44").Value = 1.5

# Reflect where the author ended up looking / selecting when they saved.
$ws.Range("C53").Select()
